# Calculate the duration of lesson plans according to the steps duration
# -> Adds a new "Sheet2" after "Sheet1" containing the lesson-plan "steps"
#    table header row (step, duration, teacher_activity, student_activity,
#    knowledge, skills, values, output, assessment_criteria, facilitator_note)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A1").Value = "step"
$ws2.Range("B1").Value = "duration"
$ws2.Range("C1").Value = "teacher_activity"
$ws2.Range("D1").Value = "student_activity"
$ws2.Range("E1").Value = "knowledge"
$ws2.Range("F1").Value = "skills"
$ws2.Range("G1").Value = "values"
$ws2.Range("H1").Value = "output"
$ws2.Range("I1").Value = "assessment_criteria"
$ws2.Range("J1").Value = "facilitator_note"

# Column widths (character units - engine rounds to its pixel grid, so we
# subtract the fixed padding offset it re-adds on save to land as close as
# possible to the authored widths)
$ws2.Columns.Item(3).ColumnWidth = 14.33203125 - (5/6)
$ws2.Columns.Item(4).ColumnWidth = 15.83203125 - (5/6)
$ws2.Columns.Item(9).ColumnWidth = 18.5 - (5/6)
$ws2.Columns.Item(10).ColumnWidth = 16 - (5/6)

# Cursor/selection position left on Sheet2 when the file was last saved
$ws2.Range("K6").Select() | Out-Null

# Sheet2 is the active sheet/tab
$ws2.Activate() | Out-Null
